# Node_relationship.xlsx — "Full Dashboard After Some Changes"
#
# The sheet's row 13 (A13="FI_1", B13="F2_5") is an exact duplicate of row 12
# (A12="FI_1", B12="F2_5"). The edit removes that duplicate row entirely,
# which shifts every following row up by one (140 data rows -> 139 data
# rows). The sheet's _FilterDatabase defined name (and the implicit used
# range) needs to shrink from B1:B140 to B1:B139 to match, and the active
# selection ends up on the full column B (as if the user had just clicked
# the column header before/after trimming the filter range).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate row; everything below (rows 14..140) shifts up to
# 13..139 automatically, shared-string usage count drops by 2 (one A-cell +
# one B-cell reference removed), and the sheet dimension becomes A1:B139.
$ws.Rows(13).Delete()

# The hidden _xlnm._FilterDatabase name still points at the old B1:B140
# range (there's no live AutoFilter to auto-shrink it), so update it by
# hand to track the new last row.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -eq "Sheet1!_FilterDatabase") {
        $nm.RefersTo = "=Sheet1!`$B`$1:`$B`$139"
    }
}

# Leave the selection on the (now one-shorter) column B, matching the
# saved workbook state.
$ws.Range("B:B").Select()
